# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# on a handful of leve rows across the workbook's sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 1994.1904
$ws.Range("I40").Value = 2839.8572
$ws.Range("J40").Value = 1571.3572
$ws.Range("K40").Value = 2839.8572
$ws.Range("L40").Value = 1571.3572
$ws.Range("M40").Value = -2664.8572
$ws.Range("N40").Value = -1921.3572

$ws.Range("H64").Value = 2850
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2850
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 2850
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3346

$ws.Range("H67").Value = 2850
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2850
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 2850
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4566

$ws.Range("H76").Value = 55320.684
$ws.Range("I76").Value = 55320.684
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 55320.684
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -55005.684
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 55320.684
$ws.Range("I79").Value = 55320.684
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 55320.684
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -54228.684
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 92700.17999999999
$ws.Range("I86").Value = 126744.125
$ws.Range("J86").Value = 1916.3334
$ws.Range("K86").Value = 126744.125
$ws.Range("L86").Value = 1916.3334
$ws.Range("M86").Value = -125621.125
$ws.Range("N86").Value = -4162.3334

$ws.Range("H89").Value = 92700.17999999999
$ws.Range("I89").Value = 126744.125
$ws.Range("J89").Value = 1916.3334
$ws.Range("K89").Value = 633720.625
$ws.Range("L89").Value = 9581.666999999999
$ws.Range("M89").Value = -628104.625
$ws.Range("N89").Value = -20813.667

$ws.Range("H132").Value = 2697544.2
$ws.Range("I132").Value = 3107642.5
$ws.Range("J132").Value = 2613.7144
$ws.Range("K132").Value = 9322927.5
$ws.Range("L132").Value = 7841.1432
$ws.Range("M132").Value = -9320397.5
$ws.Range("N132").Value = -12901.1432

$ws.Range("H138").Value = 2280.82
$ws.Range("I138").Value = 1054.1034
$ws.Range("J138").Value = 2781.8733
$ws.Range("K138").Value = 3162.3102
$ws.Range("L138").Value = 8345.619900000002
$ws.Range("M138").Value = 1977.6898
$ws.Range("N138").Value = -18625.6199

# ---- ARM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H63").Value = 7666.6665
$ws.Range("I63").Value = 6500
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 6500
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -5814
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 7666.6665
$ws.Range("I66").Value = 6500
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 32500
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -29068
$ws.Range("N66").Value = -56864

$ws.Range("H74").Value = 831.5577
$ws.Range("I74").Value = 800.8958
$ws.Range("J74").Value = 1199.5
$ws.Range("K74").Value = 800.8958
$ws.Range("L74").Value = 1199.5
$ws.Range("M74").Value = 73.10419999999999
$ws.Range("N74").Value = -2947.5

$ws.Range("H77").Value = 831.5577
$ws.Range("I77").Value = 800.8958
$ws.Range("J77").Value = 1199.5
$ws.Range("K77").Value = 4004.479
$ws.Range("L77").Value = 5997.5
$ws.Range("M77").Value = 363.5209999999997
$ws.Range("N77").Value = -14733.5

$ws.Range("H132").Value = 4870.4854
$ws.Range("I132").Value = 5534.306
$ws.Range("J132").Value = 3158.5264
$ws.Range("K132").Value = 16602.918
$ws.Range("L132").Value = 9475.5792
$ws.Range("M132").Value = -14072.918
$ws.Range("N132").Value = -14535.5792

# ---- BSM --------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 2333.2727
$ws.Range("I105").Value = 2195
$ws.Range("J105").Value = 2955.5
$ws.Range("K105").Value = 2195
$ws.Range("L105").Value = 2955.5
$ws.Range("M105").Value = -448
$ws.Range("N105").Value = -6449.5

# ---- CRP --------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2428.7234
$ws.Range("I31").Value = 1798.3158
$ws.Range("J31").Value = 5090.4443
$ws.Range("K31").Value = 1798.3158
$ws.Range("L31").Value = 5090.4443
$ws.Range("M31").Value = -1503.3158
$ws.Range("N31").Value = -5680.4443

$ws.Range("H34").Value = 2428.7234
$ws.Range("I34").Value = 1798.3158
$ws.Range("J34").Value = 5090.4443
$ws.Range("K34").Value = 1798.3158
$ws.Range("L34").Value = 5090.4443
$ws.Range("M34").Value = -1596.3158
$ws.Range("N34").Value = -5494.4443

$ws.Range("H86").Value = 62501770
$ws.Range("I86").Value = 83334910
$ws.Range("J86").Value = 2352
$ws.Range("K86").Value = 83334910
$ws.Range("L86").Value = 2352
$ws.Range("M86").Value = -83333787
$ws.Range("N86").Value = -4598

$ws.Range("H89").Value = 62501770
$ws.Range("I89").Value = 83334910
$ws.Range("J89").Value = 2352
$ws.Range("K89").Value = 416674550
$ws.Range("L89").Value = 11760
$ws.Range("M89").Value = -416668934
$ws.Range("N89").Value = -22992

$ws.Range("H107").Value = 677.05084
$ws.Range("I107").Value = 632.85
$ws.Range("J107").Value = 770.1053000000001
$ws.Range("K107").Value = 632.85
$ws.Range("L107").Value = 770.1053000000001
$ws.Range("M107").Value = 1287.15
$ws.Range("N107").Value = -4610.1053

$ws.Range("H132").Value = 2501290.5
$ws.Range("I132").Value = 957.90625
$ws.Range("J132").Value = 6946326
$ws.Range("K132").Value = 2873.71875
$ws.Range("L132").Value = 20838978
$ws.Range("M132").Value = -343.71875
$ws.Range("N132").Value = -20844038

$ws.Range("H134").Value = 2855.3333
$ws.Range("I134").Value = 2897.5
$ws.Range("K134").Value = 8692.5
$ws.Range("M134").Value = -6157.5

$ws.Range("H141").Value = 51468.727
$ws.Range("J141").Value = 58500.89
$ws.Range("L141").Value = 58500.89
$ws.Range("N141").Value = -68860.89

# ---- CUL --------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 170
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# ---- GSM --------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 16133416
$ws.Range("I70").Value = 21743492
$ws.Range("J70").Value = 4451.125
$ws.Range("K70").Value = 21743492
$ws.Range("L70").Value = 4451.125
$ws.Range("M70").Value = -21743222
$ws.Range("N70").Value = -4991.125

$ws.Range("H73").Value = 16133416
$ws.Range("I73").Value = 21743492
$ws.Range("J73").Value = 4451.125
$ws.Range("K73").Value = 21743492
$ws.Range("L73").Value = 4451.125
$ws.Range("M73").Value = -21742556
$ws.Range("N73").Value = -6323.125

$ws.Range("H80").Value = 6475
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 5000
$ws.Range("M80").Value = -4002

$ws.Range("H83").Value = 6475
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 25000
$ws.Range("M83").Value = -20008

# ---- LTW --------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H122").Value = 3450.8462
$ws.Range("I122").Value = 3857.0557
$ws.Range("J122").Value = 2536.875
$ws.Range("K122").Value = 11571.1671
$ws.Range("L122").Value = 7610.625
$ws.Range("M122").Value = -9121.167099999999
$ws.Range("N122").Value = -12510.625

$ws.Range("H136").Value = 2942.7058
$ws.Range("I136").Value = 2902.0667
$ws.Range("J136").Value = 3247.5
$ws.Range("K136").Value = 8706.2001
$ws.Range("L136").Value = 9742.5
$ws.Range("M136").Value = -6156.2001
$ws.Range("N136").Value = -14842.5

# ---- WVR --------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 1251.1143
$ws.Range("I132").Value = 1027.909
$ws.Range("J132").Value = 2069.5334
$ws.Range("K132").Value = 3083.727
$ws.Range("L132").Value = 6208.600199999999
$ws.Range("M132").Value = -553.7270000000003
$ws.Range("N132").Value = -11268.6002

$ws.Range("H136").Value = 1389.6038
$ws.Range("I136").Value = 1450.6666
$ws.Range("J136").Value = 1219.5
$ws.Range("K136").Value = 4351.9998
$ws.Range("L136").Value = 3658.5
$ws.Range("M136").Value = -1801.9998
$ws.Range("N136").Value = -8758.5
